$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2022-07-29) was added to the series.
# It belongs right above the former row 24, so insert a fresh row there
# which pushes the existing rows 24:51 down to 25:52.
$ws.Rows("24:24").Insert()

# Populate the newly inserted row 24 with the new weekly record.
$ws.Range("A24").Value = 9
$ws.Range("B24").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44771
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = 100112035
$ws.Range("G24").Value = "Bruselas (repollito)"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 43
$ws.Range("K24").Value = 22000
$ws.Range("L24").Value = 22000
$ws.Range("M24").Value = 22000
$ws.Range("N24").Value = "$/malla 15 kilos"
$ws.Range("O24").Value = "Hijuelas"
$ws.Range("P24").Value = 1467
$ws.Range("Q24").Value = 15
$ws.Range("R24").Value = "Hortaliza"
